$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.023469273611642
$ws.Range("D2").Value = 1.02991880175578
$ws.Range("E2").Value = 1.024061280717853
$ws.Range("F2").Value = 1.021935023578461
$ws.Range("I2").Value = 1.031958491711991
$ws.Range("J2").Value = 1.028649318193341
$ws.Range("K2").Value = 1.032731318345536
$ws.Range("L2").Value = 1.02689088810567
$ws.Range("M2").Value = 1.024770887594859
$ws.Range("N2").Value = 1.013524052227593
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.024663589256033
$ws.Range("D3").Value = 1.030437582726837
$ws.Range("E3").Value = 1.025082261966219
$ws.Range("F3").Value = 1.023770401479187
$ws.Range("I3").Value = 1.032173262586644
$ws.Range("J3").Value = 1.029480857179979
$ws.Range("K3").Value = 1.033059124242721
$ws.Range("L3").Value = 1.027718298054218
$ws.Range("M3").Value = 1.026410013247868
$ws.Range("N3").Value = 1.01380604068883
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.025435420836341
$ws.Range("D4").Value = 1.030773138799663
$ws.Range("E4").Value = 1.025742399570061
$ws.Range("F4").Value = 1.024956755179782
$ws.Range("I4").Value = 1.032310757899169
$ws.Range("J4").Value = 1.030017475783464
$ws.Range("K4").Value = 1.033270393957047
$ws.Range("L4").Value = 1.028252584494128
$ws.Range("M4").Value = 1.027468971693605
$ws.Range("N4").Value = 1.013987849683499
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02575967054806
$ws.Range("D5").Value = 1.030914174928463
$ws.Range("E5").Value = 1.026019803071849
$ws.Range("F5").Value = 1.025455208634882
$ws.Range("I5").Value = 1.03236820814786
$ws.Range("J5").Value = 1.030242727786418
$ws.Range("K5").Value = 1.033359010026998
$ws.Range("L5").Value = 1.028476937011484
$ws.Range("M5").Value = 1.027913769455561
$ws.Range("N5").Value = 1.014064125900913
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.025814100232109
$ws.Range("D6").Value = 1.030937853643336
$ws.Range("E6").Value = 1.026066373477965
$ws.Range("F6").Value = 1.025538884516061
$ws.Range("I6").Value = 1.032377833625436
$ws.Range("J6").Value = 1.030280528591796
$ws.Range("K6").Value = 1.033373877226245
$ws.Range("L6").Value = 1.028514591526571
$ws.Range("M6").Value = 1.027988430482393
$ws.Range("N6").Value = 1.014076923873674
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02543975436885
$ws.Range("D7").Value = 1.03077502345638
$ws.Range("E7").Value = 1.025746106710954
$ws.Range("F7").Value = 1.024963416662949
$ws.Range("I7").Value = 1.032311526937748
$ws.Range("J7").Value = 1.030020486953771
$ws.Range("K7").Value = 1.033271578841943
$ws.Range("L7").Value = 1.028255583327789
$ws.Range("M7").Value = 1.027474916616335
$ws.Range("N7").Value = 1.013988869503085
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.023873101573095
$ws.Range("D8").Value = 1.030094152788679
$ws.Range("E8").Value = 1.024406432404628
$ws.Range("F8").Value = 1.022555565752756
$ws.Range("I8").Value = 1.032031380431106
$ws.Range("J8").Value = 1.028930641184772
$ws.Range("K8").Value = 1.032842276361202
$ws.Range("L8").Value = 1.02717074514486
$ws.Range("M8").Value = 1.025325188870711
$ws.Range("N8").Value = 1.013619487978728
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02110484361987
$ws.Range("D9").Value = 1.028893399862901
$ws.Range("E9").Value = 1.022041762210305
$ws.Range("F9").Value = 1.018302462571017
$ws.Range("I9").Value = 1.031526399331679
$ws.Range("J9").Value = 1.026999011783853
$ws.Range("K9").Value = 1.032079335905584
$ws.Range("L9").Value = 1.02525055614274
$ws.Range("M9").Value = 1.021523873794377
$ws.Range("N9").Value = 1.012963524436588
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.019253959980092
$ws.Range("D10").Value = 1.028092273381889
$ws.Range("E10").Value = 1.020462458560465
$ws.Range("F10").Value = 1.015459494710596
$ws.Range("I10").Value = 1.031182094139356
$ws.Range("J10").Value = 1.025703553239431
$ws.Range("K10").Value = 1.031566365310539
$ws.Range("L10").Value = 1.023964505220054
$ws.Range("M10").Value = 1.018980123967013
$ws.Range("N10").Value = 1.012522753766633
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.018451173027334
$ws.Range("D11").Value = 1.027745234010301
$ws.Range("E11").Value = 1.019777885820017
$ws.Range("F11").Value = 1.01422650021524
$ws.Range("I11").Value = 1.03103118310848
$ws.Range("J11").Value = 1.025140735240894
$ws.Range("K11").Value = 1.03134321156159
$ws.Range("L11").Value = 1.023406187970681
$ws.Range("M11").Value = 1.01787624414125
$ws.Range("N11").Value = 1.01233106146126
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.018152775455425
$ws.Range("D12").Value = 1.027616306514151
$ws.Range("E12").Value = 1.01952349279946
$ws.Range("F12").Value = 1.013768200955014
$ws.Range("I12").Value = 1.030974853124046
$ws.Range("J12").Value = 1.024931393991291
$ws.Range("K12").Value = 1.031260166749784
$ws.Range("L12").Value = 1.023198583095058
$ws.Range("M12").Value = 1.017465837958954
$ws.Range("N12").Value = 1.01225973171622
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.018216792186138
$ws.Range("D13").Value = 1.027643962872546
$ws.Range("E13").Value = 1.019578066100487
$ws.Range("F13").Value = 1.01386652196773
$ws.Range("I13").Value = 1.03098694854976
$ws.Range("J13").Value = 1.024976311383925
$ws.Range("K13").Value = 1.0312779871874
$ws.Range("L13").Value = 1.023243125099808
$ws.Range("M13").Value = 1.017553888699041
$ws.Range("N13").Value = 1.01227503793989
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.018426511638401
$ws.Range("D14").Value = 1.027734577262304
$ws.Range("E14").Value = 1.01975685992315
$ws.Range("F14").Value = 1.014188623465719
$ws.Range("I14").Value = 1.031026532466998
$ws.Range("J14").Value = 1.025123436875589
$ws.Range("K14").Value = 1.031336350226471
$ws.Range("L14").Value = 1.023389031814453
$ws.Range("M14").Value = 1.017842327582958
$ws.Range("N14").Value = 1.012325167913279
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.018555699163973
$ws.Range("D15").Value = 1.027790404923895
$ws.Range("E15").Value = 1.019867005722864
$ws.Range("F15").Value = 1.014387039299153
$ws.Range("I15").Value = 1.031050884973274
$ws.Range("J15").Value = 1.025214047791043
$ws.Range("K15").Value = 1.031372288999478
$ws.Range("L15").Value = 1.023478900373139
$ws.Range("M15").Value = 1.01801999408861
$ws.Range("N15").Value = 1.012356037822895
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.019307209593628
$ws.Range("D16").Value = 1.02811530219051
$ws.Range("E16").Value = 1.020507875852364
$ws.Range("F16").Value = 1.015541281771046
$ws.Range("I16").Value = 1.03119207109859
$ws.Range("J16").Value = 1.025740865773826
$ws.Range("K16").Value = 1.031581153484893
$ws.Range("L16").Value = 1.024001528090282
$ws.Range("M16").Value = 1.019053332790624
$ws.Range("N16").Value = 1.012535458049967
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019778249308895
$ws.Range("D17").Value = 1.02831906274888
$ws.Range("E17").Value = 1.020909680517903
$ws.Range("F17").Value = 1.016264770034823
$ws.Range("I17").Value = 1.031280144376685
$ws.Range("J17").Value = 1.026070819905696
$ws.Range("K17").Value = 1.031711891655762
$ws.Range("L17").Value = 1.024328968441727
$ws.Range("M17").Value = 1.019700862274207
$ws.Range("N17").Value = 1.01264777908191
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.020052869658577
$ws.Range("D18").Value = 1.02843789869516
$ws.Range("E18").Value = 1.02114397677893
$ws.Range("F18").Value = 1.016686579320276
$ws.Range("I18").Value = 1.031331340035015
$ws.Range("J18").Value = 1.026263095827292
$ws.Range("K18").Value = 1.03178804925933
$ws.Range("L18").Value = 1.024519819306916
$ws.Range("M18").Value = 1.020078322688379
$ws.Range("N18").Value = 1.012713213484464
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.020146486349188
$ws.Range("D19").Value = 1.028478416271883
$ws.Range("E19").Value = 1.021223854047661
$ws.Range("F19").Value = 1.01683037372401
$ws.Range("I19").Value = 1.03134876660304
$ws.Range("J19").Value = 1.026328626403574
$ws.Range("K19").Value = 1.031814000119077
$ws.Range("L19").Value = 1.024584870949654
$ws.Range("M19").Value = 1.02020698786483
$ws.Range("N19").Value = 1.012735511295207
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.019727724601986
$ws.Range("D20").Value = 1.028297202629693
$ws.Range("E20").Value = 1.020866577905477
$ws.Range("F20").Value = 1.016187166263607
$ws.Range("I20").Value = 1.031270713157425
$ws.Range("J20").Value = 1.026035437673112
$ws.Range("K20").Value = 1.031697875007805
$ws.Range("L20").Value = 1.024293851650532
$ws.Range("M20").Value = 1.019631412638829
$ws.Range("N20").Value = 1.012635736438927
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.018364760235526
$ws.Range("D21").Value = 1.027707894177076
$ws.Range("E21").Value = 1.019704212724719
$ws.Range("F21").Value = 1.014093781280034
$ws.Range("I21").Value = 1.031014883585116
$ws.Range("J21").Value = 1.025080119999169
$ws.Range("K21").Value = 1.031319168062488
$ws.Range("L21").Value = 1.023346072050993
$ws.Range("M21").Value = 1.017757400003582
$ws.Range("N21").Value = 1.012310409391942
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.017506611666783
$ws.Range("D22").Value = 1.027337247632734
$ws.Range("E22").Value = 1.018972736560302
$ws.Range("F22").Value = 1.012775785215857
$ws.Range("I22").Value = 1.030852442449926
$ws.Range("J22").Value = 1.0244778196678
$ws.Range("K22").Value = 1.031080159759585
$ws.Range("L22").Value = 1.022748885142406
$ws.Range("M22").Value = 1.01657695178012
$ws.Range("N22").Value = 1.012105129818552
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017961647807796
$ws.Range("D23").Value = 1.027533746104633
$ws.Range("E23").Value = 1.019360568795268
$ws.Range("F23").Value = 1.013474655409771
$ws.Range("I23").Value = 1.030938706634362
$ws.Range("J23").Value = 1.024797268529326
$ws.Range("K23").Value = 1.031206947987367
$ws.Range("L23").Value = 1.023065587697629
$ws.Range("M23").Value = 1.017202940783113
$ws.Range("N23").Value = 1.012214022317829
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019750554940941
$ws.Range("D24").Value = 1.028307080320448
$ws.Range("E24").Value = 1.020886054335583
$ws.Range("F24").Value = 1.016222232650574
$ws.Range("I24").Value = 1.031274975263222
$ws.Range("J24").Value = 1.026051425939218
$ws.Range("K24").Value = 1.031704208835981
$ws.Range("L24").Value = 1.024309719847879
$ws.Range("M24").Value = 1.019662794657056
$ws.Range("N24").Value = 1.012641178239761
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.021821433581915
$ws.Range("D25").Value = 1.029203935403733
$ws.Range("E25").Value = 1.022653576940675
$ws.Range("F25").Value = 1.019403271127207
$ws.Range("I25").Value = 1.031658294983842
$ws.Range("J25").Value = 1.027499728597666
$ws.Range("K25").Value = 1.032277339561605
$ws.Range("L25").Value = 1.025748003068959
$ws.Range("M25").Value = 1.022508238273027
$ws.Range("N25").Value = 1.013133712766594
